$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same conference listing data
# and need the same updates applied (mirrors the source diff touching two
# identical worksheet XML fragments).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 4: 南宁·2024良牙动漫秋季盛典（秋典）
    $ws.Range("F4").Value = 3727
    $ws.Range("G4").Value = 58

    # Row 5: 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini
    $ws.Range("F5").Value = 384
}
